# Populate names of students and lecturers and other things.
#
# The original sheet had a "role" column (E) whose values were all the
# constant "student" - it gets removed entirely, shifting every column
# after it one to the left (ID, department, semester). A new student row
# (Eslam) is appended, a couple of existing cells are corrected, and the
# misspelled "deparment" header is fixed to "department".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused "role" column (E); everything to its right (ID,
# deparment, semester + their column-width formatting) shifts left by one.
$ws.Columns.Item(5).Delete()

# --- Fix up the rows that survive the shift -----------------------------

# nabil (row 3): department was IS, is now CS; semester was 8, is now 7.
$ws.Range("F3").Value = "CS"
$ws.Range("G3").Value = 7

# omar (row 4): level was 3, is now 4; semester was 6, is now 8.
$ws.Range("C4").Value = 4
$ws.Range("G4").Value = 8

# --- New row: Eslam -------------------------------------------------------

$ws.Range("A5").Value = "Eslam"
$ws.Range("B5").Value = "eslam@e.com"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "pass33"
$ws.Range("E5").Value = 204050
$ws.Range("F5").Value = "IS"
$ws.Range("G5").Value = 8

$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:eslam@e.com")
$ws.Range("B5").Style = "Hyperlink"

# --- Correct the "deparment" -> "department" header typo (column F now) --

$ws.Range("F1").Value = "department"

# --- Selection, matching the saved workbook state -------------------------

$ws.Range("F1").Select()
